$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timestamp in the table header (also updates the table column name)
$ws.Range("C1").Value = "2024-11-15 17:47:12"

# Update attendance status values to "Falta"
$ws.Range("C2").Value = "Falta"
$ws.Range("C3").Value = "Falta"
$ws.Range("C4").Value = "Falta"
$ws.Range("C5").Value = "Falta"
$ws.Range("C6").Value = "Falta"
$ws.Range("C7").Value = "Falta"
$ws.Range("C10").Value = "Falta"
$ws.Range("C12").Value = "Falta"
